# Insert a new slide ("Why use the pattern over a list of states?") as the
# 5th slide of the deck (pushing the former slide 5 "UML of Demo:" and
# everything after it back by one position).

$p = $ppt.ActivePresentation

# ppLayoutText (2) == the "Title, Content" AutoLayout -> title + body
# placeholder, matching the TITLE_AND_BODY slide layout used by the rest of
# this deck.
$newSlide = $p.Slides.Add(5, 2)

# --- Title -------------------------------------------------------------
$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Why use the pattern over a list of states?"
$title.ParagraphFormat.Alignment = 2   # ppAlignCenter

# --- Body ----------------------------------------------------------------
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "You don" + [char]8217 + "t have direct access to the states in a simple list."
$body.InsertAfter([char]13 + "You miss the ability to interact with the specific save states.") | Out-Null
$body.InsertAfter([char]13 + "Having the memento patterns allows for branching paths via ") | Out-Null
$body.InsertAfter("save states.") | Out-Null

Write-Output "Inserted slide at index $($newSlide.SlideIndex) (id $($newSlide.SlideID))"
